# Update gh-pages to output generated at 456a3b4
# Applies the "wanted to go" (F column) count refresh across all four
# sheets, plus a newly scraped row on 本地生活 (Local life).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------
# Sheet "展览" (Exhibitions) - refreshed "想去人数" (wanted-to-go) counts
# ---------------------------------------------------------------
$ws1.Range("F8").Value = 140
$ws1.Range("F10").Value = 12495
$ws1.Range("F11").Value = 12826
$ws1.Range("F12").Value = 1314
$ws1.Range("F14").Value = 5451
$ws1.Range("F15").Value = 0
$ws1.Range("F16").Value = 109
$ws1.Range("F18").Value = 194
$ws1.Range("F20").Value = 0
$ws1.Range("F26").Value = 503
$ws1.Range("F28").Value = 2990
$ws1.Range("F29").Value = 252
$ws1.Range("F33").Value = 0
$ws1.Range("F34").Value = 137
$ws1.Range("F37").Value = 3725
$ws1.Range("F38").Value = 0
$ws1.Range("F44").Value = 0
$ws1.Range("F45").Value = 0
$ws1.Range("F46").Value = 0
$ws1.Range("F47").Value = 33
$ws1.Range("F48").Value = 29
$ws1.Range("F49").Value = 4293

# ---------------------------------------------------------------
# Sheet "演出" (Performances) - refreshed counts
# ---------------------------------------------------------------
$ws2.Range("F4").Value = 0
$ws2.Range("F10").Value = 0
$ws2.Range("F12").Value = 0
$ws2.Range("F15").Value = 0
$ws2.Range("F20").Value = 12
$ws2.Range("F27").Value = 49
$ws2.Range("F28").Value = 2
$ws2.Range("F31").Value = 1

# ---------------------------------------------------------------
# Sheet "本地生活" (Local life) - refreshed count + new scraped row
# ---------------------------------------------------------------
$ws3.Range("F2").Value = 0

# New row 3: keep text-like columns (date, numbers-as-id) formatted the
# same way the existing row 2 is, then fall back to copying row 2's
# cell formatting onto column A/B so the new row matches existing style.
$ws3.Range("A3").Value = 2
$ws3.Range("B3").NumberFormat = "@"
$ws3.Range("B3").Value = "2024-07-14"
$ws3.Range("B3").Style = "Normal"
$ws3.Range("C3").Value = "杭州·排球少年!!垃圾场决战 主题咖啡厅"
$ws3.Range("D3").Value = "延安路292号（地铁1号线龙翔桥站D出口） 工联CC"
$ws3.Range("E3").Value = "2024.07.14 00:00-08.25 23:59"
$ws3.Range("F3").Value = 6
$ws3.Range("G3").Value = 10
$ws3.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=89031"
$ws3.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202407/v0pLPU0J1720599572561.jpeg"

# Match A2's look (bold/centered/bordered "index" style) on the new A3 cell.
$ws3.Range("A2").Copy()
$ws3.Range("A3").PasteSpecial(-4122)
$ws3.Range("A3").Value = 2

# ---------------------------------------------------------------
# Sheet "全部类型" (All types) - refreshed counts
# ---------------------------------------------------------------
$ws4.Range("F2").Value = 570
$ws4.Range("F4").Value = 234
$ws4.Range("F6").Value = 6978
$ws4.Range("F7").Value = 0
$ws4.Range("F8").Value = 0
$ws4.Range("F11").Value = 0
$ws4.Range("F18").Value = 0
$ws4.Range("F20").Value = 351
$ws4.Range("F21").Value = 2011
$ws4.Range("F24").Value = 0
$ws4.Range("F26").Value = 731
$ws4.Range("F29").Value = 0
$ws4.Range("F31").Value = 6
$ws4.Range("F33").Value = 6
$ws4.Range("F36").Value = 137
$ws4.Range("F40").Value = 0
$ws4.Range("F44").Value = 644
$ws4.Range("F45").Value = 1258
$ws4.Range("F46").Value = 913
$ws4.Range("F48").Value = 29
$ws4.Range("F49").Value = 0
$ws4.Range("F50").Value = 0

Write-Output "done"
